$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 20 by shifting row 19's data down, then rewrite rows 16-19
# with their updated values, as described by the diff.

# First, capture row 19 values (the row that will become row 20)
$row19Vals = @{}
for ($col = 1; $col -le 20; $col++) {
    $row19Vals[$col] = $ws.Cells.Item(19, $col).Value
}

# Update dimension / extend used range by writing to row 20 (new last row)
# Row 20 = old row 19 data (unchanged values), since that whole record moved down
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44910
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100101
$ws.Cells.Item(20, 8).Value = "Berries"
$ws.Cells.Item(20, 9).Value = 100101004
$ws.Cells.Item(20, 10).Value = "Frambuesa"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 200
$ws.Cells.Item(20, 14).Value = 7500
$ws.Cells.Item(20, 15).Value = 8000
$ws.Cells.Item(20, 16).Value = 7750
$ws.Cells.Item(20, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(20, 18).Value = "Región de Ñuble"
$ws.Cells.Item(20, 19).Value = 3875
$ws.Cells.Item(20, 20).Value = 2

# Row 19 now takes on values that used to be associated with a "Segunda" record
$ws.Cells.Item(19, 4).Value = 44532
$ws.Cells.Item(19, 12).Value = "Segunda"
$ws.Cells.Item(19, 13).Value = 100
$ws.Cells.Item(19, 14).Value = 8000
$ws.Cells.Item(19, 15).Value = 8000
$ws.Cells.Item(19, 16).Value = 8000
$ws.Cells.Item(19, 19).Value = 4000

# Row 18 updates
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 14).Value = 10000
$ws.Cells.Item(18, 15).Value = 10000
$ws.Cells.Item(18, 16).Value = 10000
$ws.Cells.Item(18, 19).Value = 5000

# Row 17 updates
$ws.Cells.Item(17, 4).Value = 44944
$ws.Cells.Item(17, 14).Value = 7000
$ws.Cells.Item(17, 15).Value = 8000
$ws.Cells.Item(17, 16).Value = 7500
$ws.Cells.Item(17, 19).Value = 3750

# Row 16 updates
$ws.Cells.Item(16, 4).Value = 44988
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 6000
$ws.Cells.Item(16, 15).Value = 7000
$ws.Cells.Item(16, 16).Value = 6500
$ws.Cells.Item(16, 19).Value = 3250
